$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (previously row 8's data)
$ws.Range("D2").Value = 44232
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11583
$ws.Range("S2").Value = 827

# Update row 4 (previously row 5's data)
$ws.Range("D4").Value = 44229
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 11000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 11364
$ws.Range("S4").Value = 812

# Update row 5 (previously row 7's data)
$ws.Range("D5").Value = 44210
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 11000
$ws.Range("P5").Value = 10357
$ws.Range("S5").Value = 740

# Update row 6 (previously row 4's data)
$ws.Range("D6").Value = 44253
$ws.Range("M6").Value = 90
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 12667
$ws.Range("S6").Value = 905

# Update row 7 (previously row 2's data)
$ws.Range("D7").Value = 44216
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11545
$ws.Range("S7").Value = 825

# Update row 8 (previously row 6's data)
$ws.Range("D8").Value = 44181
$ws.Range("M8").Value = 65
$ws.Range("N8").Value = 9000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 9462
$ws.Range("S8").Value = 676
